$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update activation date (keep it stored as plain text, just like the
#    original "01/01/2018" value, instead of letting Excel auto-convert a
#    date-look-alike literal into a date serial number / new number
#    format style). Building the text via a formula first, then freezing
#    the computed result back onto the cell as a plain value (copy /
#    paste-special-values onto itself) avoids Excel's literal-entry date
#    autodetection entirely, while keeping the existing cell style.
$ws.Range("B8").Formula = "=""01/01/"" & ""2022"""
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Formula = "=""01/01/"" & ""2022"""
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# 2) Add new English objectives text in row 11 (B11 / C11), matching the
#    styling pattern used throughout column B (s=2) and column C (s=3).
#    Copy the formatting from the row above (same column styles) instead
#    of toggling individual font/alignment properties, which avoids
#    generating spurious extra cell-style entries in styles.xml.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newObjectivesText = "Provide students with the knowledge of cell biology necessary to understand the other subjects of the course and the training of the Environmental Engineer."
$ws.Range("B11").Value = $newObjectivesText
$ws.Range("C11").Value = $newObjectivesText

# 3) Update "Programa resumido" texts (row 14): remove the leading
#    "Origem e evolução da célula; " phrase.
$newResumo = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."
$ws.Range("B14").Value = $newResumo
$ws.Range("C14").Value = $newResumo

# 4) Update "Short syllabus" texts (row 15): remove the leading
#    "The origin and evolution of the cell; " phrase.
$newShortSyllabus = "Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division."
$ws.Range("B15").Value = $newShortSyllabus
$ws.Range("C15").Value = $newShortSyllabus

# 5) Update "Programa" texts (row 16): replace opening bullet about the
#    origin/evolution of the cell with the new wording about cell
#    structure and evolutionary history.
$newPrograma = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."
$ws.Range("B16").Value = $newPrograma
$ws.Range("C16").Value = $newPrograma

# 6) Update "Syllabus" texts (row 17): same change, in English.
$newSyllabus = "Cell structure and evolutionary history: prokaryotic microorganisms andeukaryotic and their evolutionary relationships between the Bacteria, Archaea andEukarya.Microscope analysis of cells structure: optical and electron microscope.Structure and function of major organic molecules: carbohydrates, lipids, nucleic acids and proteins. Internal organization of the cell: membrane structure and function; intracelular compartments and protein sorting; vesicular traffic (endocytosis and exocytosis).Nucleus and genetic material organization: structure and functionCell cycle and cell division: mitosis and meiosisCell energy conversion: mitochondria and chloroplast."
$ws.Range("B17").Value = $newSyllabus
$ws.Range("C17").Value = $newSyllabus
